# Rename "Requested quantity" headers on existing sheets
$wb = $excel.ActiveWorkbook
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")

$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# Add the new "PO Forecast" sheet after the last existing sheet
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Header row
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Copy header formatting (bold, border, centered) from the Weekly Quantity sheet
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

$data = @(
    @(45151.99999999999, 22, 4.401891069727118, 42.74816373845292),
    @(45158.99999999999, 23, 1.786645325414139, 42.12642068421562),
    @(45172.99999999999, 24, 3.141141755291179, 42.81373740250532),
    @(45200.99999999999, 26, 6.887598727978665, 45.02657536466641),
    @(45242.99999999999, 30, 9.165626000288203, 49.45308884619443),
    @(45256.99999999999, 31, 11.25099737799092, 50.43545952143466),
    @(45277.99999999999, 32, 12.95494666994289, 51.54441380489734),
    @(45298.99999999999, 34, 14.24161855989757, 51.97777751105575),
    @(45312.99999999999, 35, 15.82977225531105, 54.2975666374078),
    @(45326.99999999999, 37, 17.96068290060986, 55.25281483061599),
    @(45368.99999999999, 40, 20.07410268762527, 59.52211476717285),
    @(45424.99999999999, 45, 25.14979583610558, 63.2439591806989),
    @(45438.99999999999, 46, 25.99174365568077, 63.83426684822549),
    @(45480.99999999999, 49, 31.14046868021089, 69.73707450858528),
    @(45536.99999999999, 54, 34.05067829507013, 74.52998252025851),
    @(45564.99999999999, 56, 36.46349987971987, 76.30143982863638),
    @(45578.99999999999, 57, 37.79426829001281, 75.49334833535777),
    @(45585.99999999999, 58, 39.59647713991943, 77.7707478612122),
    @(45592.99999999999, 59, 39.12650259993435, 78.51331102878335),
    @(45599.99999999999, 59, 38.95218476135896, 78.30783374023494),
    @(45606.99999999999, 60, 40.45405823266355, 78.30045340939442),
    @(45613.99999999999, 60, 40.5962238114616, 78.29190465943039),
    @(45620.99999999999, 61, 40.54525880409545, 79.20805234824651),
    @(45627.99999999999, 62, 43.88636064083975, 79.93102935078653),
    @(45634.99999999999, 62, 42.72502737705257, 81.74820429430106),
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $wsForecast.Cells.Item($row, 1).Value = $data[$i][0]
    $wsForecast.Cells.Item($row, 2).Value = $data[$i][1]
    $wsForecast.Cells.Item($row, 3).Value = $data[$i][2]
    $wsForecast.Cells.Item($row, 4).Value = $data[$i][3]
}

# Copy the date-column formatting (numeric date style) from Weekly Quantity A2 to the new A2:A26 range
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A26").PasteSpecial(-4122)
